$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows: fill in the "Sr.No." column (A2:A4) ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Row 4 gained content in F4 (reuses existing "Anything Relevant" string)
$ws.Range("F4").Value = "Anything Relevant"

# New shared strings must be introduced in this exact order so they
# land at the same sharedStrings.xml indices as the target workbook:
#   15 = amazon, 16 = a, 17 = TEST-002
$ws.Range("F6").Value = "amazon"
$ws.Range("G4").Value = "a"
$ws.Range("B5").Value = "TEST-002"

# Row 5 (TEST-002 test case, mirrors TEST-001 in row 2)
$ws.Range("A5").Value = 4
$ws.Range("C5").Value = "Validate whether user can search with any valid phrase"
$ws.Range("D5").Value = "Y"
$ws.Range("E5").Value = "go_to_google_home"

$ws.Range("A6").Value = 5
$ws.Range("E6").Value = "search_for_phrase"

$ws.Range("A7").Value = 6
$ws.Range("E7").Value = "validate_first_link_text"
$ws.Range("F7").Value = "Anything Relevant"
$ws.Range("G7").Value = "a"

# Copy formatting (styles) from the template rows into the new rows
$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A3:G3").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)

$ws.Range("A4:G4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)

# Match the row heights Excel computed when it wrapped the new text
$ws.Rows(5).RowHeight = 31
$ws.Rows(6).RowHeight = 15.5
$ws.Rows(7).RowHeight = 46.5

$excel.CutCopyMode = $false

$ws.Range("D5").Select()
